$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B36 should be a real number (2) instead of an inline string "2"
$ws.Range("B36").Value = 2

# Add new row 37 with annotation data
$ws.Range("A37").Value = "Sunsi Wu"
$ws.Range("B37").Value = "'3"
$ws.Range("C37").Value = "无"
$ws.Range("D37").Value = "DIS"
$ws.Range("E37").Value = "MET"
$ws.Range("F37").Value = "0c8a854c-e7df-48dd-93a0-b6771319a745"
$ws.Range("G37").Value = "H1Ww66x0-_annotated.xlsx"
$ws.Range("H37").Value = "- the proposed approach to maintain the budget is simplistic"
